# Phase 4.4 metadata refresh: Experimental flag flips to false and the
# generation Date timestamp is bumped. Both values live in the "Metadata"
# worksheet (Property/Value pairs) as plain text ("true"/"false" are stored
# as strings in this workbook, not native booleans), so a straightforward
# Range.Value assignment must be avoided: Excel auto-coerces the literal
# strings "true"/"false" into boolean cells (t="b"), which does not match
# the source data. Instead we write the text through a formula that
# evaluates to the desired string, then collapse the formula down to a
# static value with a values-only paste so the cell keeps its original
# text type and style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# Experimental: true -> false
Set-TextValue $ws.Range("B7") "false"

# Date: refreshed generation timestamp
Set-TextValue $ws.Range("B8") "2025-10-03T16:37:46+01:00"

$excel.CutCopyMode = 0
